$d = $word.ActiveDocument
$count = 0

$r = $d.Content
$found = $r.Find.Execute("13×28=364", $true, $false, $false, $false, $false, $true, 1, $false, "26×59=1534", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 13×28=364" }

$r = $d.Content
$found = $r.Find.Execute("42×78=3276", $true, $false, $false, $false, $false, $true, 1, $false, "36×10=360", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 42×78=3276" }

$r = $d.Content
$found = $r.Find.Execute("18×90=1620", $true, $false, $false, $false, $false, $true, 1, $false, "93×13=1209", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 18×90=1620" }

$r = $d.Content
$found = $r.Find.Execute("54×43=2322", $true, $false, $false, $false, $false, $true, 1, $false, "16×66=1056", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 54×43=2322" }

$r = $d.Content
$found = $r.Find.Execute("16×44=704", $true, $false, $false, $false, $false, $true, 1, $false, "65×24=1560", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 16×44=704" }

$r = $d.Content
$found = $r.Find.Execute("83×45=3735", $true, $false, $false, $false, $false, $true, 1, $false, "72×76=5472", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 83×45=3735" }

$r = $d.Content
$found = $r.Find.Execute("31×28=868", $true, $false, $false, $false, $false, $true, 1, $false, "44×37=1628", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 31×28=868" }

$r = $d.Content
$found = $r.Find.Execute("35×15=525", $true, $false, $false, $false, $false, $true, 1, $false, "58×74=4292", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 35×15=525" }

$r = $d.Content
$found = $r.Find.Execute("23×18=414", $true, $false, $false, $false, $false, $true, 1, $false, "14×68=952", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 23×18=414" }

$r = $d.Content
$found = $r.Find.Execute("66×24=1584", $true, $false, $false, $false, $false, $true, 1, $false, "60×76=4560", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 66×24=1584" }

$r = $d.Content
$found = $r.Find.Execute("71×87=6177", $true, $false, $false, $false, $false, $true, 1, $false, "31×14=434", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 71×87=6177" }

$r = $d.Content
$found = $r.Find.Execute("25×26=650", $true, $false, $false, $false, $false, $true, 1, $false, "19×17=323", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 25×26=650" }

$r = $d.Content
$found = $r.Find.Execute("46×57=2622", $true, $false, $false, $false, $false, $true, 1, $false, "63×83=5229", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 46×57=2622" }

$r = $d.Content
$found = $r.Find.Execute("40×12=480", $true, $false, $false, $false, $false, $true, 1, $false, "35×11=385", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 40×12=480" }

$r = $d.Content
$found = $r.Find.Execute("93×20=1860", $true, $false, $false, $false, $false, $true, 1, $false, "41×62=2542", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 93×20=1860" }

$r = $d.Content
$found = $r.Find.Execute("99×10=990", $true, $false, $false, $false, $false, $true, 1, $false, "47×28=1316", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 99×10=990" }

$r = $d.Content
$found = $r.Find.Execute("45×32=1440", $true, $false, $false, $false, $false, $true, 1, $false, "11×54=594", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 45×32=1440" }

$r = $d.Content
$found = $r.Find.Execute("79×63=4977", $true, $false, $false, $false, $false, $true, 1, $false, "86×80=6880", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 79×63=4977" }

$r = $d.Content
$found = $r.Find.Execute("100×49=4900", $true, $false, $false, $false, $false, $true, 1, $false, "30×80=2400", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 100×49=4900" }

$r = $d.Content
$found = $r.Find.Execute("90×57=5130", $true, $false, $false, $false, $false, $true, 1, $false, "82×13=1066", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 90×57=5130" }

$r = $d.Content
$found = $r.Find.Execute("58×12=696", $true, $false, $false, $false, $false, $true, 1, $false, "36×94=3384", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 58×12=696" }

$r = $d.Content
$found = $r.Find.Execute("31×48=1488", $true, $false, $false, $false, $false, $true, 1, $false, "85×92=7820", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 31×48=1488" }

$r = $d.Content
$found = $r.Find.Execute("13×35=455", $true, $false, $false, $false, $false, $true, 1, $false, "78×61=4758", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 13×35=455" }

$r = $d.Content
$found = $r.Find.Execute("31×51=1581", $true, $false, $false, $false, $false, $true, 1, $false, "98×54=5292", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 31×51=1581" }

$r = $d.Content
$found = $r.Find.Execute("25×86=2150", $true, $false, $false, $false, $false, $true, 1, $false, "13×39=507", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 25×86=2150" }

$r = $d.Content
$found = $r.Find.Execute("12×31=372", $true, $false, $false, $false, $false, $true, 1, $false, "29×76=2204", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 12×31=372" }

$r = $d.Content
$found = $r.Find.Execute("46×15=690", $true, $false, $false, $false, $false, $true, 1, $false, "27×74=1998", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 46×15=690" }

$r = $d.Content
$found = $r.Find.Execute("16×92=1472", $true, $false, $false, $false, $false, $true, 1, $false, "17×29=493", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 16×92=1472" }

$r = $d.Content
$found = $r.Find.Execute("59×59=3481", $true, $false, $false, $false, $false, $true, 1, $false, "36×77=2772", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 59×59=3481" }

$r = $d.Content
$found = $r.Find.Execute("57×93=5301", $true, $false, $false, $false, $false, $true, 1, $false, "79×60=4740", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 57×93=5301" }

$r = $d.Content
$found = $r.Find.Execute("76×59=4484", $true, $false, $false, $false, $false, $true, 1, $false, "23×79=1817", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 76×59=4484" }

$r = $d.Content
$found = $r.Find.Execute("83×17=1411", $true, $false, $false, $false, $false, $true, 1, $false, "26×81=2106", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 83×17=1411" }

$r = $d.Content
$found = $r.Find.Execute("16×46=736", $true, $false, $false, $false, $false, $true, 1, $false, "20×49=980", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 16×46=736" }

$r = $d.Content
$found = $r.Find.Execute("57×94=5358", $true, $false, $false, $false, $false, $true, 1, $false, "17×18=306", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 57×94=5358" }

$r = $d.Content
$found = $r.Find.Execute("84×42=3528", $true, $false, $false, $false, $false, $true, 1, $false, "13×49=637", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 84×42=3528" }

$r = $d.Content
$found = $r.Find.Execute("63×15=945", $true, $false, $false, $false, $false, $true, 1, $false, "33×95=3135", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 63×15=945" }

$r = $d.Content
$found = $r.Find.Execute("59×46=2714", $true, $false, $false, $false, $false, $true, 1, $false, "51×74=3774", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 59×46=2714" }

$r = $d.Content
$found = $r.Find.Execute("44×47=2068", $true, $false, $false, $false, $false, $true, 1, $false, "19×94=1786", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 44×47=2068" }

$r = $d.Content
$found = $r.Find.Execute("22×42=924", $true, $false, $false, $false, $false, $true, 1, $false, "58×29=1682", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 22×42=924" }

$r = $d.Content
$found = $r.Find.Execute("29×64=1856", $true, $false, $false, $false, $false, $true, 1, $false, "94×17=1598", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 29×64=1856" }

$r = $d.Content
$found = $r.Find.Execute("57×36=2052", $true, $false, $false, $false, $false, $true, 1, $false, "26×91=2366", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 57×36=2052" }

$r = $d.Content
$found = $r.Find.Execute("69×27=1863", $true, $false, $false, $false, $false, $true, 1, $false, "10×97=970", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 69×27=1863" }

$r = $d.Content
$found = $r.Find.Execute("77×51=3927", $true, $false, $false, $false, $false, $true, 1, $false, "66×17=1122", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 77×51=3927" }

$r = $d.Content
$found = $r.Find.Execute("43×12=516", $true, $false, $false, $false, $false, $true, 1, $false, "56×94=5264", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 43×12=516" }

$r = $d.Content
$found = $r.Find.Execute("70×36=2520", $true, $false, $false, $false, $false, $true, 1, $false, "53×21=1113", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 70×36=2520" }

$r = $d.Content
$found = $r.Find.Execute("50×98=4900", $true, $false, $false, $false, $false, $true, 1, $false, "61×69=4209", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 50×98=4900" }

$r = $d.Content
$found = $r.Find.Execute("96×21=2016", $true, $false, $false, $false, $false, $true, 1, $false, "44×13=572", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 96×21=2016" }

$r = $d.Content
$found = $r.Find.Execute("96×93=8928", $true, $false, $false, $false, $false, $true, 1, $false, "31×43=1333", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 96×93=8928" }

$r = $d.Content
$found = $r.Find.Execute("52×38=1976", $true, $false, $false, $false, $false, $true, 1, $false, "82×27=2214", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 52×38=1976" }

$r = $d.Content
$found = $r.Find.Execute("72×46=3312", $true, $false, $false, $false, $false, $true, 1, $false, "56×66=3696", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 72×46=3312" }

$r = $d.Content
$found = $r.Find.Execute("74×73=5402", $true, $false, $false, $false, $false, $true, 1, $false, "47×12=564", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 74×73=5402" }

$r = $d.Content
$found = $r.Find.Execute("88×38=3344", $true, $false, $false, $false, $false, $true, 1, $false, "90×45=4050", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 88×38=3344" }

$r = $d.Content
$found = $r.Find.Execute("49×11=539", $true, $false, $false, $false, $false, $true, 1, $false, "40×21=840", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 49×11=539" }

$r = $d.Content
$found = $r.Find.Execute("16×24=384", $true, $false, $false, $false, $false, $true, 1, $false, "80×31=2480", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 16×24=384" }

$r = $d.Content
$found = $r.Find.Execute("99×38=3762", $true, $false, $false, $false, $false, $true, 1, $false, "68×25=1700", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 99×38=3762" }

$r = $d.Content
$found = $r.Find.Execute("17×99=1683", $true, $false, $false, $false, $false, $true, 1, $false, "91×57=5187", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 17×99=1683" }

$r = $d.Content
$found = $r.Find.Execute("68×57=3876", $true, $false, $false, $false, $false, $true, 1, $false, "96×96=9216", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 68×57=3876" }

$r = $d.Content
$found = $r.Find.Execute("16×15=240", $true, $false, $false, $false, $false, $true, 1, $false, "53×48=2544", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 16×15=240" }

$r = $d.Content
$found = $r.Find.Execute("14×14=196", $true, $false, $false, $false, $false, $true, 1, $false, "84×39=3276", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 14×14=196" }

$r = $d.Content
$found = $r.Find.Execute("54×89=4806", $true, $false, $false, $false, $false, $true, 1, $false, "82×81=6642", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 54×89=4806" }

$r = $d.Content
$found = $r.Find.Execute("50×48=2400", $true, $false, $false, $false, $false, $true, 1, $false, "83×78=6474", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 50×48=2400" }

$r = $d.Content
$found = $r.Find.Execute("66×18=1188", $true, $false, $false, $false, $false, $true, 1, $false, "74×66=4884", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 66×18=1188" }

$r = $d.Content
$found = $r.Find.Execute("21×99=2079", $true, $false, $false, $false, $false, $true, 1, $false, "45×56=2520", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 21×99=2079" }

$r = $d.Content
$found = $r.Find.Execute("85×72=6120", $true, $false, $false, $false, $false, $true, 1, $false, "73×40=2920", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 85×72=6120" }

$r = $d.Content
$found = $r.Find.Execute("17×15=255", $true, $false, $false, $false, $false, $true, 1, $false, "99×25=2475", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 17×15=255" }

$r = $d.Content
$found = $r.Find.Execute("51×39=1989", $true, $false, $false, $false, $false, $true, 1, $false, "89×36=3204", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 51×39=1989" }

$r = $d.Content
$found = $r.Find.Execute("28×56=1568", $true, $false, $false, $false, $false, $true, 1, $false, "24×98=2352", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 28×56=1568" }

$r = $d.Content
$found = $r.Find.Execute("74×81=5994", $true, $false, $false, $false, $false, $true, 1, $false, "85×97=8245", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 74×81=5994" }

$r = $d.Content
$found = $r.Find.Execute("73×55=4015", $true, $false, $false, $false, $false, $true, 1, $false, "35×85=2975", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 73×55=4015" }

$r = $d.Content
$found = $r.Find.Execute("66×31=2046", $true, $false, $false, $false, $false, $true, 1, $false, "34×50=1700", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 66×31=2046" }

$r = $d.Content
$found = $r.Find.Execute("68×45=3060", $true, $false, $false, $false, $false, $true, 1, $false, "15×19=285", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 68×45=3060" }

$r = $d.Content
$found = $r.Find.Execute("65×47=3055", $true, $false, $false, $false, $false, $true, 1, $false, "61×74=4514", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 65×47=3055" }

$r = $d.Content
$found = $r.Find.Execute("67×47=3149", $true, $false, $false, $false, $false, $true, 1, $false, "78×37=2886", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 67×47=3149" }

$r = $d.Content
$found = $r.Find.Execute("66×99=6534", $true, $false, $false, $false, $false, $true, 1, $false, "60×42=2520", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 66×99=6534" }

$r = $d.Content
$found = $r.Find.Execute("82×72=5904", $true, $false, $false, $false, $false, $true, 1, $false, "75×86=6450", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 82×72=5904" }

$r = $d.Content
$found = $r.Find.Execute("94×49=4606", $true, $false, $false, $false, $false, $true, 1, $false, "80×66=5280", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 94×49=4606" }

$r = $d.Content
$found = $r.Find.Execute("86×27=2322", $true, $false, $false, $false, $false, $true, 1, $false, "48×86=4128", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 86×27=2322" }

$r = $d.Content
$found = $r.Find.Execute("80×97=7760", $true, $false, $false, $false, $false, $true, 1, $false, "30×28=840", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 80×97=7760" }

$r = $d.Content
$found = $r.Find.Execute("13×56=728", $true, $false, $false, $false, $false, $true, 1, $false, "25×78=1950", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 13×56=728" }

$r = $d.Content
$found = $r.Find.Execute("47×63=2961", $true, $false, $false, $false, $false, $true, 1, $false, "39×31=1209", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 47×63=2961" }

$r = $d.Content
$found = $r.Find.Execute("76×77=5852", $true, $false, $false, $false, $false, $true, 1, $false, "36×72=2592", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 76×77=5852" }

$r = $d.Content
$found = $r.Find.Execute("10×82=820", $true, $false, $false, $false, $false, $true, 1, $false, "60×83=4980", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 10×82=820" }

$r = $d.Content
$found = $r.Find.Execute("80×64=5120", $true, $false, $false, $false, $false, $true, 1, $false, "17×46=782", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 80×64=5120" }

$r = $d.Content
$found = $r.Find.Execute("95×89=8455", $true, $false, $false, $false, $false, $true, 1, $false, "50×89=4450", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 95×89=8455" }

$r = $d.Content
$found = $r.Find.Execute("55×87=4785", $true, $false, $false, $false, $false, $true, 1, $false, "34×37=1258", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 55×87=4785" }

$r = $d.Content
$found = $r.Find.Execute("78×18=1404", $true, $false, $false, $false, $false, $true, 1, $false, "62×75=4650", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 78×18=1404" }

$r = $d.Content
$found = $r.Find.Execute("67×23=1541", $true, $false, $false, $false, $false, $true, 1, $false, "81×11=891", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 67×23=1541" }

$r = $d.Content
$found = $r.Find.Execute("50×69=3450", $true, $false, $false, $false, $false, $true, 1, $false, "55×41=2255", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 50×69=3450" }

$r = $d.Content
$found = $r.Find.Execute("98×31=3038", $true, $false, $false, $false, $false, $true, 1, $false, "49×39=1911", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 98×31=3038" }

$r = $d.Content
$found = $r.Find.Execute("87×10=870", $true, $false, $false, $false, $false, $true, 1, $false, "53×100=5300", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 87×10=870" }

$r = $d.Content
$found = $r.Find.Execute("20×20=400", $true, $false, $false, $false, $false, $true, 1, $false, "10×51=510", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 20×20=400" }

$r = $d.Content
$found = $r.Find.Execute("63×37=2331", $true, $false, $false, $false, $false, $true, 1, $false, "24×17=408", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 63×37=2331" }

$r = $d.Content
$found = $r.Find.Execute("90×95=8550", $true, $false, $false, $false, $false, $true, 1, $false, "33×20=660", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 90×95=8550" }

$r = $d.Content
$found = $r.Find.Execute("65×57=3705", $true, $false, $false, $false, $false, $true, 1, $false, "87×63=5481", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 65×57=3705" }

$r = $d.Content
$found = $r.Find.Execute("14×90=1260", $true, $false, $false, $false, $false, $true, 1, $false, "62×26=1612", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 14×90=1260" }

$r = $d.Content
$found = $r.Find.Execute("15×36=540", $true, $false, $false, $false, $false, $true, 1, $false, "35×62=2170", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 15×36=540" }

$r = $d.Content
$found = $r.Find.Execute("75×40=3000", $true, $false, $false, $false, $false, $true, 1, $false, "23×72=1656", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 75×40=3000" }

$r = $d.Content
$found = $r.Find.Execute("99×23=2277", $true, $false, $false, $false, $false, $true, 1, $false, "55×89=4895", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 99×23=2277" }

$r = $d.Content
$found = $r.Find.Execute("28×68=1904", $true, $false, $false, $false, $false, $true, 1, $false, "16×85=1360", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 28×68=1904" }

$r = $d.Content
$found = $r.Find.Execute("32×35=1120", $true, $false, $false, $false, $false, $true, 1, $false, "97×54=5238", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 32×35=1120" }

Write-Output "Replaced $count of 100"